$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (totals) sheet: insert a new 2022-Q3 row at the top of the data
#    (row 2), pushing the existing 2022-Q2 / 2021-Q4 / 2021-Q2 rows down by
#    one. Easiest/most robust way: write the new last row (row 5) first by
#    cloning the format of row 4, then overwrite rows 2-4 with the values
#    that slide down from the old rows 2-3-4.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Clone the formatting of the last existing data row (row 4, style "s=2" on
# column A) onto the new row 5 so the new row matches the sheet's existing
# look, then fill in row 5 with what used to be row 4's data (2021-Q2).
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)  # xlPasteFormats

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 1.19

# Row 4 becomes what used to be row 3 (2021-Q4)
$total.Range("B4").Value = "2021-Q4"
$total.Range("C4").Value = 5
$total.Range("D4").Value = 0.44

# Row 3 becomes what used to be row 2 (2022-Q2)
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.58

# Row 2 becomes the brand-new 2022-Q3 entry
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.29

# ---------------------------------------------------------------------------
# 2) Duplicate the "2022-Q2" sheet (same fund lineup reports a new quarter),
#    place the copy right before it, rename to "2022-Q3" and overwrite the
#    quarter-specific figures with the new values.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Row 2: 华夏磐利一年定期开放混合A (009686)
$q3.Range("D2").Value = "'10.76"
$q3.Range("E2").Value = "'64.78"
$q3.Range("F2").Value = "'2.57"
$q3.Range("G2").Value = "'0.2765"
$q3.Range("H2").Value = 6

# Row 3: 华夏磐利一年定期开放混合C (009687)
$q3.Range("D3").Value = "'0.43"
$q3.Range("E3").Value = "'64.78"
$q3.Range("F3").Value = "'2.57"
$q3.Range("G3").Value = "'0.0111"
$q3.Range("H3").Value = 6

Write-Output "2022-Q3 sheet added"
